$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update/clarify the definitions ("Definition" / "Display definition" columns)
# of four recently-added terms (EUPATH_0000138..0000141).
$ws.Range("E10").Value = "A time datum for the duration of a patients participation in an investigation. "
$ws.Range("F11").Value = "average number of clinical visits where data was collected for a patient, average over calendar years"
$ws.Range("F12").Value = "average number of clinical visits where data was collected for patient and the patient was diagnosed with asymptomatic parasitemia, average over calendar years"
$ws.Range("F13").Value = "average number of clinical visits where data was collected for patient nd the patient was diagnosed with malaria, average over calendar years"

# Trim the large block of unused trailing blank/placeholder rows so the
# worksheet's used range shrinks from A1:J52 down to A1:J21, keeping only a
# small buffer of blank formatted rows below the data table.
$ws.Rows("14:38").Delete()
$ws.Rows("22:27").Delete()

# Move the active selection back to the top of the data table.
$ws.Range("A2").Select() | Out-Null
